# "Query and table name updated"
# - Rename the worksheet from "Students" to "STUDENTS"
# - Rename the "First Name"/"Last Name" columns (header cells + table
#   column definitions) to "First_Name"/"Last_Name"
# - Leave active cell on B1 (matches the post-edit selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Table/worksheet name change
$ws.Name = "STUDENTS"

# Header + table column renames (ListObject header cells drive the
# table's column names automatically when the source range is a table)
$ws.Range("A1").Value = "First_Name"
$ws.Range("B1").Value = "Last_Name"

# Final selection left on B1
$ws.Range("B1").Select()
